$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: row, open, close, high, low, shares_outstanding
$data = @(
    @(2, 32.59999847412109, 32.41999816894531, 34.93000030517578, 31.78000068664551, 31367897),
    @(3, 32.11000061035156, 30.79000091552734, 32.88999938964844, 30.10000038146973, 31367897),
    @(4, 29.06999969482422, 30.01000022888184, 32.13000106811523, 28.27000045776367, 31367897),
    @(5, 34.34000015258789, 33.93000030517578, 34.9900016784668, 30, 31367897),
    @(6, 35.02000045776367, 37.38000106811523, 38.54000091552734, 34.15000152587891, 31367897),
    @(7, 34.9900016784668, 37.91999816894531, 38.40000152587891, 34.27000045776367, 31367897),
    @(8, 43.7400016784668, 50.70000076293945, 50.84000015258789, 42.06000137329102, 31367897),
    @(9, 59.22000122070312, 60.56000137329102, 63.16999816894531, 57.34000015258789, 31367897),
    @(10, 66.11000061035156, 79.73999786376953, 82.44000244140625, 65, 31367897),
    @(11, 83.30000305175781, 105.5100021362305, 109.7799987792969, 79.4000015258789, 31367897),
    @(12, 77, 67.19000244140625, 81.93000030517578, 62.95999908447266, 31367897),
    @(13, 77.80999755859375, 88.75, 94.8499984741211, 77.68000030517578, 31367897),
    @(14, 63.97999954223633, 58.36000061035156, 68.12000274658203, 56.90000152587891, 31367897),
    @(15, 75.0999984741211, 79.25, 89.4800033569336, 70.86000061035156, 31367897),
    @(16, 111.3300018310547, 88.16999816894531, 114.2200012207031, 76.11000061035156, 31367897),
    @(17, 82.05999755859375, 89.79000091552734, 92.47000122070312, 77.9800033569336, 31367897),
    @(18, 118.3099975585938, 103.5400009155273, 132.3000030517578, 96.5, 31367897),
    @(19, 93.9499969482422, 73.0999984741211, 96.0999984741211, 69.25, 31367897),
    @(20, 61.77000045776367, 48.06000137329102, 62.38999938964844, 44, 31367897),
    @(21, 51.29999923706055, 47.04999923706055, 57.90000152587891, 46.68999862670898, 31367897),
    @(22, 39.56999969482422, 60.95999908447266, 62.13000106811523, 37.0099983215332, 31367897),
    @(23, 84, 72.69000244140625, 89.26000213623047, 70.76000213623047, 31367897),
    @(24, 76.75, 79.58000183105469, 98.5500030517578, 75.4000015258789, 31367897),
    @(25, 80.47000122070312, 120.0999984741211, 197.3999938964844, 78.54000091552734, 31367897),
    @(26, 122.870002746582, 108.8000030517578, 125.3499984741211, 107.9199981689453, 31367897),
    @(27, 93.62999725341795, 87.5, 94.4800033569336, 80.55999755859375, 31367897),
    @(28, 79.08999633789062, 83.41999816894531, 85.41999816894531, 76.20999908447266, 31367897),
    @(29, 66.33999633789062, 65.5199966430664, 70.30000305175781, 58.43999862670898, 31367897),
    @(30, 63.88000106811523, 50.65000152587891, 65.61000061035156, 48.36000061035156, 31367897),
    @(31, 37.13000106811523, 46.0099983215332, 46.29999923706055, 36.34999847412109, 31367897),
    @(32, 56.47999954223633, 56.5, 56.84999847412109, 55.34999847412109, 31367897),
    @(33, 48.22000122070312, 45, 48.81000137329102, 44.97999954223633, 31367897),
    @(34, 43.5, 39.33000183105469, 43.7400016784668, 39.2400016784668, 31367897),
    @(35, 45.25, 40, 51.4900016784668, 38.90000152587891, 31367897),
    @(36, 37.90000152587891, 32.93000030517578, 39.02999877929688, 32.65000152587891, 31367897),
    @(37, 38.5, 13.60000038146973, 38.86999893188477, 13.5, 31367897),
    @(38, 8.800000190734863, 8.5600004196167, 10.76000022888184, 6.480000019073486, 31367897),
    @(39, 9.050000190734863, 11.78999996185303, 13.44999980926514, 8.359999656677246, 31367897),
    @(40, 8.6899995803833, 8.729999542236328, 8.920000076293945, 7.690000057220459, 31367897),
    @(41, 7.880000114440918, 7.690000057220459, 13.0600004196167, 7.449999809265137, 31367897),
    @(42, 2.690000057220459, 2.450000047683716, 2.690000057220459, 1.759999990463257, 31367897),
    @(43, 3.150000095367432, 4.210000038146973, 5.519999980926514, 3.049999952316284, 31367897)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value  = $item[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value  = $item[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value  = $item[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value  = $item[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value  = $item[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value  = "IRBT"     # I: fixed_ticker
}
